# Update register.xlsx: refresh ticket register rows 2-4 with the new
# receipt numbers / timestamps / durations / totals, and drop the old
# rows 5-9 so the register only keeps the 3 current entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ticket 310001 -> 317000)
$ws.Range("A2").Value = "'317000"
$ws.Range("C2").Value = "22/07/2024 15:43:00"
$ws.Range("D2").Value = "22/07/2024 18:31:00"
$ws.Range("G2").Value = 10080
$ws.Range("H2").Value = 4500

# Row 3 (ticket 310002 -> 317001)
$ws.Range("A3").Value = "'317001"
$ws.Range("C3").Value = "22/07/2024 19:14:00"
$ws.Range("D3").Value = "22/07/2024 19:18:00"
$ws.Range("G3").Value = 240
$ws.Range("H3").Value = 1500

# Row 4 (ticket 310003 -> 317002)
$ws.Range("A4").Value = "'317002"
$ws.Range("C4").Value = "22/07/2024 19:17:00"
$ws.Range("D4").Value = "22/07/2024 19:25:00"
$ws.Range("G4").Value = 480
$ws.Range("H4").Value = 1500

# Remove old rows 5-9 (tickets 310004-310008 no longer in the register)
$ws.Range("A5:H9").EntireRow.Delete()
